# Update the date heading paragraph
$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)
if ($p1.Range.Text.TrimEnd([char]13, [char]7) -eq "2025-06-29 Sunday") {
    $p1.Range.Text = "2025-06-30 Monday"
} else {
    throw "Unexpected date paragraph text: " + $p1.Range.Text
}

# Update the practice-problem answers in the first table, cell by cell,
# addressed by (row, column) so the duplicate text "83÷4=20, 3" (which
# appears both as an old value and as a new value elsewhere) is handled
# unambiguously.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "29÷2=14, 1") {
    $cell.Range.Text = "61÷8=7, 5"
} else {
    throw "Unexpected cell (1,1) text: " + $cell.Range.Text
}

$cell = $t.Cell(1, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "93÷7=13, 2") {
    $cell.Range.Text = "65÷3=21, 2"
} else {
    throw "Unexpected cell (1,2) text: " + $cell.Range.Text
}

$cell = $t.Cell(1, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "14÷5=2, 4") {
    $cell.Range.Text = "35÷2=17, 1"
} else {
    throw "Unexpected cell (1,3) text: " + $cell.Range.Text
}

$cell = $t.Cell(1, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "30÷7=4, 2") {
    $cell.Range.Text = "93÷9=10, 3"
} else {
    throw "Unexpected cell (1,4) text: " + $cell.Range.Text
}

$cell = $t.Cell(1, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "33÷2=16, 1") {
    $cell.Range.Text = "64÷2=32, 0"
} else {
    throw "Unexpected cell (1,5) text: " + $cell.Range.Text
}

$cell = $t.Cell(5, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "23÷7=3, 2") {
    $cell.Range.Text = "83÷4=20, 3"
} else {
    throw "Unexpected cell (5,1) text: " + $cell.Range.Text
}

$cell = $t.Cell(5, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "44÷4=11, 0") {
    $cell.Range.Text = "51÷4=12, 3"
} else {
    throw "Unexpected cell (5,2) text: " + $cell.Range.Text
}

$cell = $t.Cell(5, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "31÷9=3, 4") {
    $cell.Range.Text = "84÷9=9, 3"
} else {
    throw "Unexpected cell (5,3) text: " + $cell.Range.Text
}

$cell = $t.Cell(5, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "27÷3=9, 0") {
    $cell.Range.Text = "26÷4=6, 2"
} else {
    throw "Unexpected cell (5,4) text: " + $cell.Range.Text
}

$cell = $t.Cell(5, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "10÷2=5, 0") {
    $cell.Range.Text = "81÷7=11, 4"
} else {
    throw "Unexpected cell (5,5) text: " + $cell.Range.Text
}

$cell = $t.Cell(9, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "20÷3=6, 2") {
    $cell.Range.Text = "92÷8=11, 4"
} else {
    throw "Unexpected cell (9,1) text: " + $cell.Range.Text
}

$cell = $t.Cell(9, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "78÷7=11, 1") {
    $cell.Range.Text = "92÷6=15, 2"
} else {
    throw "Unexpected cell (9,2) text: " + $cell.Range.Text
}

$cell = $t.Cell(9, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "40÷8=5, 0") {
    $cell.Range.Text = "63÷3=21, 0"
} else {
    throw "Unexpected cell (9,3) text: " + $cell.Range.Text
}

$cell = $t.Cell(9, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "87÷3=29, 0") {
    $cell.Range.Text = "47÷6=7, 5"
} else {
    throw "Unexpected cell (9,4) text: " + $cell.Range.Text
}

$cell = $t.Cell(9, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "84÷2=42, 0") {
    $cell.Range.Text = "92÷8=11, 4"
} else {
    throw "Unexpected cell (9,5) text: " + $cell.Range.Text
}

$cell = $t.Cell(13, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "95÷3=31, 2") {
    $cell.Range.Text = "79÷4=19, 3"
} else {
    throw "Unexpected cell (13,1) text: " + $cell.Range.Text
}

$cell = $t.Cell(13, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "28÷9=3, 1") {
    $cell.Range.Text = "71÷9=7, 8"
} else {
    throw "Unexpected cell (13,2) text: " + $cell.Range.Text
}

$cell = $t.Cell(13, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "43÷6=7, 1") {
    $cell.Range.Text = "90÷5=18, 0"
} else {
    throw "Unexpected cell (13,3) text: " + $cell.Range.Text
}

$cell = $t.Cell(13, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "56÷2=28, 0") {
    $cell.Range.Text = "32÷2=16, 0"
} else {
    throw "Unexpected cell (13,4) text: " + $cell.Range.Text
}

$cell = $t.Cell(13, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "83÷4=20, 3") {
    $cell.Range.Text = "81÷5=16, 1"
} else {
    throw "Unexpected cell (13,5) text: " + $cell.Range.Text
}

$cell = $t.Cell(17, 1)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "25÷7=3, 4") {
    $cell.Range.Text = "38÷8=4, 6"
} else {
    throw "Unexpected cell (17,1) text: " + $cell.Range.Text
}

$cell = $t.Cell(17, 2)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "54÷3=18, 0") {
    $cell.Range.Text = "53÷6=8, 5"
} else {
    throw "Unexpected cell (17,2) text: " + $cell.Range.Text
}

$cell = $t.Cell(17, 3)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "60÷8=7, 4") {
    $cell.Range.Text = "46÷7=6, 4"
} else {
    throw "Unexpected cell (17,3) text: " + $cell.Range.Text
}

$cell = $t.Cell(17, 4)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "42÷2=21, 0") {
    $cell.Range.Text = "52÷7=7, 3"
} else {
    throw "Unexpected cell (17,4) text: " + $cell.Range.Text
}

$cell = $t.Cell(17, 5)
if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "21÷2=10, 1") {
    $cell.Range.Text = "40÷2=20, 0"
} else {
    throw "Unexpected cell (17,5) text: " + $cell.Range.Text
}
